$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a handful of existing daily totals (rows 2, 6, 15, 17 — July/2025 entries)
$ws.Range("B2").Value = 18050.18
$ws.Range("B6").Value = 24062.21
$ws.Range("B15").Value = 2474.05
$ws.Range("B17").Value = 8172.16

# Insert a new daily record for day 23 (07/2025) right after the existing row 17,
# shifting every subsequent row down by one.
$ws.Rows.Item(18).Insert()

$ws.Range("A18").Value = 23
$ws.Range("B18").Value = 11311.63
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 2025
$ws.Range("E18").Value = "07/2025"
